$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text

    if ($text -like "*{{date}}*") {
        # Date line: drop the right-justification.
        $p.Range.ParagraphFormat.Alignment = 0
    }
    elseif ($text -like "*{{subject}}*") {
        # Subject line: bold lives on both the paragraph mark (pPr/rPr)
        # and the run - clearing Bold on the whole paragraph range
        # (mark included) clears both.
        $p.Range.Font.Bold = 0
    }
    elseif (($text -like "*{{recipient_name}}*") -or ($text -like "*{{sender_name}}*")) {
        # These runs are bold but the paragraph mark itself is not -
        # use a format-only Find/Replace over the paragraph so only the
        # bold run text is affected, leaving the paragraph mark alone.
        $r = $p.Range
        $r.Find.ClearFormatting()
        $r.Find.Replacement.ClearFormatting()
        $r.Find.Text = ""
        $r.Find.Font.Bold = $true
        $r.Find.Replacement.Text = ""
        $r.Find.Replacement.Font.Bold = $false
        [void]$r.Find.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2)
    }
}
